$p = $ppt.ActivePresentation

# --- Slide 2: "Projektidee / Ziele" ---
# Move "T." from the start of the "T.Abplanalp's" run to the end of the
# preceding run, i.e.
#   "Kleinteilelager Buchungssystem für " / "T.Abplanalp's"
# becomes
#   "Kleinteilelager Buchungssystem für T. " / "Abplanalp's"
$s2 = $p.Slides.Item(2)
$sh2 = $s2.Shapes.Item(2)
$para2 = $sh2.TextFrame.TextRange.Paragraphs(1, 1)

# Update the 2nd run's characters first (range offsets stay valid since we
# go from the back of the paragraph towards the front).
$run2 = $para2.Characters(36, 13)
$run2.Text = "Abplanalp’s"

$run1 = $para2.Characters(1, 35)
$run1.Text = "Kleinteilelager Buchungssystem für T. "

# --- Slide 4: "Client" ---
# Prefix the "Javascript Klasse" heading with "Tabellen als ".
$s4 = $p.Slides.Item(4)
$sh4 = $s4.Shapes.Item(2)
$para4 = $sh4.TextFrame.TextRange.Paragraphs(1, 1)
$para4.InsertBefore("Tabellen als ")

# --- Remove the "Demo" slide (slide 5) ---
$p.Slides.Item(5).Delete()
